# Update auto scs (lamda_1 / lamda_2) and time-in-ms table, and trim the
# auto-capacity dictionary rows that no longer exist (rows 51-53 removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..50 (A, B, C, D, E)
$data = @(
    @(0, 33.94444444444444, 1.95, 0, 0.134),
    @(1, 33.94444444444444, 1.95, 3, 0.003),
    @(2, 33.94444444444444, 1.95, 4, 0.007),
    @(3, 33.94444444444444, 1.95, 5, 0.017),
    @(4, 33.94444444444444, 1.95, 6, 0.03),
    @(5, 33.94444444444444, 1.95, 7, 0.054),
    @(6, 33.94444444444444, 1.95, 8, 0.061),
    @(7, 33.94444444444444, 1.95, 9, 0.038),
    @(8, 33.94444444444444, 1.95, 10, 0.035),
    @(9, 33.94444444444444, 1.95, 11, 0.029),
    @(10, 33.94444444444444, 1.95, 12, 0.023),
    @(11, 33.94444444444444, 1.95, 13, 0.031),
    @(12, 33.94444444444444, 1.95, 14, 0.028),
    @(13, 33.94444444444444, 1.95, 15, 0.024),
    @(14, 33.94444444444444, 1.95, 16, 0.042),
    @(15, 33.94444444444444, 1.95, 17, 0.051),
    @(16, 33.94444444444444, 1.95, 18, 0.035),
    @(17, 33.94444444444444, 1.95, 19, 0.022),
    @(18, 33.94444444444444, 1.95, 20, 0.035),
    @(19, 33.94444444444444, 1.95, 21, 0.026),
    @(20, 33.94444444444444, 1.95, 22, 0.021),
    @(21, 33.94444444444444, 1.95, 23, 0.017),
    @(22, 33.94444444444444, 1.95, 24, 0.02),
    @(23, 33.94444444444444, 1.95, 25, 0.021),
    @(24, 33.94444444444444, 1.95, 26, 0.014),
    @(25, 33.94444444444444, 1.95, 27, 0.023),
    @(26, 33.94444444444444, 1.95, 28, 0.017),
    @(27, 33.94444444444444, 1.95, 29, 0.016),
    @(28, 33.94444444444444, 1.95, 30, 0.019),
    @(29, 33.94444444444444, 1.95, 31, 0.012),
    @(30, 33.94444444444444, 1.95, 32, 0.014),
    @(31, 33.94444444444444, 1.95, 33, 0.012),
    @(32, 33.94444444444444, 1.95, 34, 0.011),
    @(33, 33.94444444444444, 1.95, 35, 0.005),
    @(34, 33.94444444444444, 1.95, 36, 0.007),
    @(35, 33.94444444444444, 1.95, 37, 0.002),
    @(36, 33.94444444444444, 1.95, 38, 0.009000000000000001),
    @(37, 33.94444444444444, 1.95, 39, 0.007),
    @(38, 33.94444444444444, 1.95, 40, 0.006),
    @(39, 33.94444444444444, 1.95, 41, 0.004),
    @(40, 33.94444444444444, 1.95, 42, 0.002),
    @(41, 33.94444444444444, 1.95, 44, 0.002),
    @(42, 33.94444444444444, 1.95, 45, 0.002),
    @(43, 33.94444444444444, 1.95, 46, 0.003),
    @(44, 33.94444444444444, 1.95, 48, 0.001),
    @(45, 33.94444444444444, 1.95, 49, 0.002),
    @(46, 33.94444444444444, 1.95, 51, 0.001),
    @(47, 33.94444444444444, 1.95, 55, 0.003),
    @(48, 33.94444444444444, 1.95, 56, 0.001)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Rows 51-53 no longer exist in the refreshed dictionary - remove them so
# the sheet's used range shrinks back down to A1:E50.
$ws.Range("A51:E53").EntireRow.Delete()
